$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 55581820
$ws.Range("I62").Value = 111113150
$ws.Range("K62").Value = 111113150
$ws.Range("M62").Value = -111112526
$ws.Range("H65").Value = 55581820
$ws.Range("I65").Value = 111113150
$ws.Range("K65").Value = 555565750
$ws.Range("M65").Value = -555562630
$ws.Range("H98").Value = 6418.125
$ws.Range("I98").Value = 6508.864
$ws.Range("J98").Value = 5420
$ws.Range("K98").Value = 6508.864
$ws.Range("L98").Value = 5420
$ws.Range("M98").Value = -5010.864
$ws.Range("N98").Value = -8416
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 2339.5833
$ws.Range("I100").Value = 1461.8334
$ws.Range("K100").Value = 1461.8334
$ws.Range("M100").Value = -920.8334
$ws.Range("H103").Value = 392.0625
$ws.Range("I103").Value = 196.97144
$ws.Range("K103").Value = 590.91432
$ws.Range("M103").Value = -4.914319999999975
$ws.Range("H112").Value = 5246.1333
$ws.Range("J112").Value = 5421.6514
$ws.Range("L112").Value = 16264.9542
$ws.Range("N112").Value = -18480.9542
$ws.Range("H113").Value = 47015720
$ws.Range("J113").Value = 62510252
$ws.Range("L113").Value = 62510252
$ws.Range("N113").Value = -62516760
$ws.Range("H116").Value = 41678916
$ws.Range("I116").Value = 125004250
$ws.Range("J116").Value = 16247.5
$ws.Range("K116").Value = 125004250
$ws.Range("L116").Value = 16247.5
$ws.Range("M116").Value = -125000808
$ws.Range("N116").Value = -23131.5
$ws.Range("H122").Value = 6418.125
$ws.Range("I122").Value = 6508.864
$ws.Range("J122").Value = 5420
$ws.Range("K122").Value = 19526.592
$ws.Range("L122").Value = 16260
$ws.Range("M122").Value = -17076.592
$ws.Range("N122").Value = -21160
$ws.Range("H132").Value = 2169.6943
$ws.Range("I132").Value = 1600.6333
$ws.Range("K132").Value = 4801.8999
$ws.Range("M132").Value = -2271.8999
$ws.Range("H137").Value = 2407.8096
$ws.Range("I137").Value = 2747.4167
$ws.Range("K137").Value = 8242.250100000001
$ws.Range("M137").Value = -5692.250100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3207.5
$ws.Range("I2").Value = 2230.5454
$ws.Range("J2").Value = 3923.9333
$ws.Range("K2").Value = 2230.5454
$ws.Range("L2").Value = 3923.9333
$ws.Range("M2").Value = -2117.5454
$ws.Range("N2").Value = -4149.933300000001
$ws.Range("H116").Value = 3207.5
$ws.Range("I116").Value = 2230.5454
$ws.Range("J116").Value = 3923.9333
$ws.Range("K116").Value = 2230.5454
$ws.Range("L116").Value = 3923.9333
$ws.Range("M116").Value = 63.45460000000003
$ws.Range("N116").Value = -8511.933300000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3207.5
$ws.Range("I3").Value = 2230.5454
$ws.Range("J3").Value = 3923.9333
$ws.Range("K3").Value = 2230.5454
$ws.Range("L3").Value = 3923.9333
$ws.Range("M3").Value = -2116.5454
$ws.Range("N3").Value = -4151.933300000001
$ws.Range("H107").Value = 40184116
$ws.Range("I107").Value = 56254904
$ws.Range("K107").Value = 56254904
$ws.Range("M107").Value = -56252984
$ws.Range("H134").Value = 6965.0312
$ws.Range("I134").Value = 2631.8823
$ws.Range("K134").Value = 7895.646900000001
$ws.Range("M134").Value = -5360.646900000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4920.36
$ws.Range("I16").Value = 3210
$ws.Range("J16").Value = 6499.154
$ws.Range("K16").Value = 3210
$ws.Range("L16").Value = 6499.154
$ws.Range("M16").Value = -2923
$ws.Range("N16").Value = -7073.154
$ws.Range("H31").Value = 8598.521000000001
$ws.Range("I31").Value = 3766
$ws.Range("J31").Value = 12315.846
$ws.Range("K31").Value = 3766
$ws.Range("L31").Value = 12315.846
$ws.Range("M31").Value = -3471
$ws.Range("N31").Value = -12905.846
$ws.Range("H34").Value = 8598.521000000001
$ws.Range("I34").Value = 3766
$ws.Range("J34").Value = 12315.846
$ws.Range("K34").Value = 3766
$ws.Range("L34").Value = 12315.846
$ws.Range("M34").Value = -3564
$ws.Range("N34").Value = -12719.846
$ws.Range("H58").Value = 16136503
$ws.Range("I58").Value = 45455468
$ws.Range("J58").Value = 11071.9
$ws.Range("K58").Value = 45455468
$ws.Range("L58").Value = 11071.9
$ws.Range("M58").Value = -45455265
$ws.Range("N58").Value = -11477.9
$ws.Range("H113").Value = 4920.36
$ws.Range("I113").Value = 3210
$ws.Range("J113").Value = 6499.154
$ws.Range("K113").Value = 3210
$ws.Range("L113").Value = 6499.154
$ws.Range("M113").Value = -1040
$ws.Range("N113").Value = -10839.154
$ws.Range("H134").Value = 9095.32
$ws.Range("I134").Value = 4654.222
$ws.Range("J134").Value = 11593.4375
$ws.Range("K134").Value = 13962.666
$ws.Range("L134").Value = 34780.3125
$ws.Range("M134").Value = -11427.666
$ws.Range("N134").Value = -39850.3125
$ws.Range("H136").Value = 16136503
$ws.Range("I136").Value = 45455468
$ws.Range("J136").Value = 11071.9
$ws.Range("K136").Value = 136366404
$ws.Range("L136").Value = 33215.7
$ws.Range("M136").Value = -136363854
$ws.Range("N136").Value = -38315.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 134709.33
$ws.Range("I2").Value = 62.8
$ws.Range("K2").Value = 376.8
$ws.Range("M2").Value = -263.8
$ws.Range("H5").Value = 3640173
$ws.Range("J5").Value = 5557
$ws.Range("L5").Value = 16671
$ws.Range("N5").Value = -16895
$ws.Range("H115").Value = 1637.25
$ws.Range("I115").Value = 1274.5
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 3823.5
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = -2648.5
$ws.Range("N115").Value = -8350
$ws.Range("H121").Value = 20001300
$ws.Range("I121").Value = 100000000
$ws.Range("J121").Value = 1625.75
$ws.Range("K121").Value = 300000000
$ws.Range("L121").Value = 4877.25
$ws.Range("M121").Value = -299998690
$ws.Range("N121").Value = -7497.25
$ws.Range("H132").Value = 10235.24
$ws.Range("I132").Value = 4530.615
$ws.Range("K132").Value = 40775.535
$ws.Range("M132").Value = -38245.535
$ws.Range("H135").Value = 3640173
$ws.Range("J135").Value = 5557
$ws.Range("L135").Value = 50013
$ws.Range("N135").Value = -55083
$ws.Range("H140").Value = 167834.83
$ws.Range("I140").Value = 167834.83
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 503504.49
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -498324.49
$ws.Range("N140").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7264.484
$ws.Range("I113").Value = 3599.889
$ws.Range("J113").Value = 8763.637000000001
$ws.Range("K113").Value = 3599.889
$ws.Range("L113").Value = 8763.637000000001
$ws.Range("M113").Value = -1429.889
$ws.Range("N113").Value = -13103.637
$ws.Range("H122").Value = 1728037.6
$ws.Range("I122").Value = 2500872
$ws.Range("J122").Value = 4022.2307
$ws.Range("K122").Value = 7502616
$ws.Range("L122").Value = 12066.6921
$ws.Range("M122").Value = -7500166
$ws.Range("N122").Value = -16966.6921

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4336.778
$ws.Range("I61").Value = 1200.2142
$ws.Range("K61").Value = 1200.2142
$ws.Range("M61").Value = -998.2141999999999
$ws.Range("H113").Value = 4336.778
$ws.Range("I113").Value = 1200.2142
$ws.Range("K113").Value = 1200.2142
$ws.Range("M113").Value = 969.7858000000001
$ws.Range("H132").Value = 7941814
$ws.Range("I132").Value = 13515857
$ws.Range("K132").Value = 40547571
$ws.Range("M132").Value = -40545041

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6454897.5
$ws.Range("I81").Value = 1384.5714
$ws.Range("J81").Value = 11769555
$ws.Range("K81").Value = 2769.1428
$ws.Range("L81").Value = 23539110
$ws.Range("M81").Value = -1708.1428
$ws.Range("N81").Value = -23541232
$ws.Range("H84").Value = 6454897.5
$ws.Range("I84").Value = 1384.5714
$ws.Range("J84").Value = 11769555
$ws.Range("K84").Value = 13845.714
$ws.Range("L84").Value = 117695550
$ws.Range("M84").Value = -8541.714
$ws.Range("N84").Value = -117706158
$ws.Range("H107").Value = 15152358
$ws.Range("I107").Value = 678.5625
$ws.Range("J107").Value = 55556836
$ws.Range("K107").Value = 2035.6875
$ws.Range("L107").Value = 166670508
$ws.Range("M107").Value = -115.6875
$ws.Range("N107").Value = -166674348
$ws.Range("H124").Value = 51598
$ws.Range("J124").Value = 51598
$ws.Range("L124").Value = 51598
$ws.Range("N124").Value = -61418
$ws.Range("H136").Value = 31254536
$ws.Range("I136").Value = 62500572
$ws.Range("K136").Value = 187501716
$ws.Range("M136").Value = -187499166

Write-Host "Applied all profit sheet updates"